# Expanded number of virtual machines to 32 and updated stack depths
#
# Adds a new worksheet ("Sheet1") at the end of the workbook that documents
# and calculates the virtual-machine count together with the various stack
# / user-area depths, and makes it the active sheet/tab (the previously
# active "Virtualization" tab becomes unselected).

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# --- Text labels (entered in this order so new shared strings line up) ---
$ws.Range("H1").Value = "width"
$ws.Range("I1").Value = "depth"

$ws.Range("A2").Value = "vmp_w "
$ws.Range("A3").Value = "psp_w "
$ws.Range("A4").Value = "rsp_w "
$ws.Range("A5").Value = "ssp_w "
$ws.Range("A6").Value = "esp_w "

$ws.Range("D3").Value = "Parameter stack depth"
$ws.Range("D4").Value = "Return stack depth"
$ws.Range("D5").Value = "Subroutine stack depth"
$ws.Range("D6").Value = "Exception stack depth"

$ws.Range("D7").Value = "User data area"
$ws.Range("A7").Value = "user_w"
$ws.Range("D2").Value = "Number of virtual machines"

$ws.Range("F3").Value = "cells"
$ws.Range("F4").Value = "cells"
$ws.Range("F5").Value = "cells"
$ws.Range("F6").Value = "cells"
$ws.Range("F7").Value = "longwords"
$ws.Range("F2").Value = "instances"

# --- Numbers -------------------------------------------------------------
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 7
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 9

$ws.Range("H3").Value = 32
$ws.Range("H4").Value = 32
$ws.Range("H5").Value = 544
$ws.Range("H6").Value = 304
$ws.Range("H7").Value = 32

# --- Formulas (ranges set together become shared formulas, like Excel) ---
$ws.Range("E2").Formula = "=2^B2"
$ws.Range("E3:E6").Formula = "=2^B3"
$ws.Range("E7").Formula = "=2^B7"

$ws.Range("I3").Formula = "=2^(`$B`$2+B3)"
$ws.Range("I4:I7").Formula = "=2^(`$B`$2+B4)"

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 26.917
$ws.Range("E1:F1").ColumnWidth = 14.584
$ws.Range("I1:J1").ColumnWidth = 9.584

# --- Alignment / indentation formatting -----------------------------------
$ws.Range("E1:F1").HorizontalAlignment = -4152
$ws.Range("H1:I1").HorizontalAlignment = -4152

$ws.Range("E2:E7").HorizontalAlignment = -4152
$ws.Range("E2:E7").IndentLevel = 1

$ws.Range("F2:F7").HorizontalAlignment = -4131

$ws.Range("B5").HorizontalAlignment = -4152

$ws.Range("A1").Select()
